$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.399.48'
$ws.Range("E2").Value = '  +4.09%  '

$ws.Range("D3").Value = '1.724.84'
$ws.Range("E3").Value = '  +3.88%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9981'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.94%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4729'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2639'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.22%  '

$ws.Range("D10").Value = '1.716.91'
$ws.Range("E10").Value = '  +3.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07084'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.86%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.5952'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("E14").Value = '  +0.66%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.85%  '

$ws.Range("E16").Value = '  -0.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9992'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.07%  '

$ws.Range("D18").Value = '26.378.82'

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006820'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.80%  '

$ws.Range("E20").Value = '  +1.76%  '

$ws.Range("D21").Value = '1.937.40'
$ws.Range("E21").Value = '  +3.49%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.564'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.53%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.787'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.346'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.16%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.27'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.411'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.60%  '

$ws.Range("B28").Value = 'BitcoinCash'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '108.39'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.24%  '

$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.772'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.029'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.84%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.705'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.93%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.07767'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04472'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.612'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9797'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6225'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '116.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +19.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9287'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.39%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.420'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.89%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.918'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.81%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01481'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.376'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +15.94%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3834'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1167'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.293'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.76%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05291'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.82%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.698'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.16%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.3399'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.222'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.50%  '
